$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cell that previously contained "Assert" should now read "Then",
# to match BDD syntax (Given/When/Then).
$ws.Range("A12").Value = "Then"

# Update the conditional formatting rule that highlighted the "Assert"
# keyword so it now highlights "Then" instead.
foreach ($fc in $ws.Range("A1:XFD1048576").FormatConditions) {
    if ($fc.Operator -eq 3 -and $fc.Formula1 -eq '="Assert"') {
        $fc.Formula1 = '="Then"'
    }
}
